# Against human mode, conditional spy placing added weight
#
# Adds two new boolean columns (P: IsPlaceSpyForEnemySpy, Q: IsPlaceSpyForEnemyTroops)
# to the "Cards" sheet, populates them for every data row, extends the
# AutoFilter / _FilterDatabase range to cover the new columns, and nudges
# the view state (selection) to match the authored workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells -------------------------------------------------
$ws.Range("P1").Value = "IsPlaceSpyForEnemySpy"
$ws.Range("Q1").Value = "IsPlaceSpyForEnemyTroops"

# --- Default every data row (2-126) to FALSE / FALSE -------------------
for ($r = 2; $r -le 126; $r++) {
    $ws.Cells.Item($r, 16).Value = $false
    $ws.Cells.Item($r, 17).Value = $false
}

# --- Cards that flip IsPlaceSpyForEnemySpy (column P) to TRUE ----------
$pTrueRows = @(49)
foreach ($r in $pTrueRows) {
    $ws.Cells.Item($r, 16).Value = $true
}

# --- Cards that flip IsPlaceSpyForEnemyTroops (column Q) to TRUE -------
$qTrueRows = @(8, 28, 42, 52, 65, 84, 93, 117)
foreach ($r in $qTrueRows) {
    $ws.Cells.Item($r, 17).Value = $true
}

# --- Extend the filter range to include the new columns ----------------
$ws.Range("A1:Q126").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range
$filterDbName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterDbName.RefersTo = "=Cards!`$A`$1:`$Q`$126"

# --- View tweaks: keep header frozen, move selection -------------------
$excel.ActiveWindow.ScrollRow = 88
$ws.Range("A4:H5").Select()
